$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46043
$ws.Range("B2").Value = 87.26000000000001
$ws.Range("C2").Value = 79.59999999999999
$ws.Range("D2").Value = 73.66
$ws.Range("E2").Value = 65.87
$ws.Range("F2").Value = 64.63
$ws.Range("G2").Value = 69.65000000000001
$ws.Range("H2").Value = 84.53
$ws.Range("I2").Value = 92.11
$ws.Range("J2").Value = 97
$ws.Range("K2").Value = 94.38
$ws.Range("L2").Value = 88.66
$ws.Range("M2").Value = 80.89
$ws.Range("N2").Value = 70.61
$ws.Range("O2").Value = 69.52
$ws.Range("P2").Value = 74.88
$ws.Range("Q2").Value = 80.77
$ws.Range("R2").Value = 81
$ws.Range("S2").Value = 81.90000000000001
$ws.Range("T2").Value = 80.90000000000001
$ws.Range("U2").Value = 92.61
$ws.Range("V2").Value = 99.22
$ws.Range("W2").Value = 95.31999999999999
$ws.Range("X2").Value = 93.12
$ws.Range("Y2").Value = 79.02
$ws.Range("Z2").Value = 82.38
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 91.67
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 97.27
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 95.69
